$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'63.069.17"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  -2.05%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'3.142.58"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  +0.15%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  -0.02%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'587.50"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -2.59%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'137.22"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  -4.50%  "
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +0.05%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'3.133.97"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +0.01%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.515"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -1.52%  "
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  -3.35%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'5.22"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -2.99%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.457"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -2.80%  "
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -4.34%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'34.13"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -3.03%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'3.654.32"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -0.16%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +0.92%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'3.132.75"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  -0.27%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'63.026.78"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -2.37%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'6.64"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -3.51%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'470.36"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -2.40%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'14.13"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -3.06%  "
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -1.98%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'7.64"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -0.71%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'84.72"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -2.81%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'12.95"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  -3.56%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +0.02%  "
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -2.20%  "
$c.Style = "Normal"
$c = $ws.Range("B28")
$c.Value = "'NEARProtocol"
$c.Style = "Normal"
$c = $ws.Range("C28")
$c.Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'6.98"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -3.08%  "
$c.Style = "Normal"
$c = $ws.Range("B29")
$c.Value = "'RenderToken"
$c.Style = "Normal"
$c = $ws.Range("C29")
$c.Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'7.93"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -4.90%  "
$c.Style = "Normal"
$c = $ws.Range("B30")
$c.Value = "'ImmutableX"
$c.Style = "Normal"
$c = $ws.Range("C30")
$c.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'2.11"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +2.32%  "
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -0.11%  "
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -0.34%  "
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  -5.53%  "
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  -6.09%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'1.07"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -2.84%  "
$c.Style = "Normal"
$c = $ws.Range("B36")
$c.Value = "'Filecoin"
$c.Style = "Normal"
$c = $ws.Range("C36")
$c.Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'5.77"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  -3.99%  "
$c.Style = "Normal"
$c = $ws.Range("B37")
$c.Value = "'OKB"
$c.Style = "Normal"
$c = $ws.Range("C37")
$c.Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'52.23"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.77%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'0.0" + ([char]0x2083) + "0690"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -2.03%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'417.10"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -6.24%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'2.72"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -8.71%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'8.17"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -0.91%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'2.890.91"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +0.62%  "
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  -6.54%  "
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -0.23%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'2.11"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -5.70%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'25.36"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -2.70%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.113"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -0.91%  "
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  -8.34%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'120.47"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -1.00%  "
$c.Style = "Normal"
